$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.334.21"
$ws.Range("E2").Value = "  -4.09%  "
$ws.Range("D3").Value = "2.974.41"
$ws.Range("E3").Value = "  -6.05%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.47%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.974.03"
$ws.Range("E8").Value = "  -5.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("E10").Value = "  -6.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.10"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("E13").Value = "  -5.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.39"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.14%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "3.465.48"
$ws.Range("E16").Value = "  -5.99%  "
$ws.Range("D17").Value = "60.300.60"
$ws.Range("E17").Value = "  -4.13%  "
$ws.Range("D18").Value = "2.975.08"
$ws.Range("E18").Value = "  -6.21%  "
$ws.Range("E19").Value = "  -6.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "423.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.13%  "
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.81"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.35%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -6.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.76%  "
$ws.Range("E30").Value = "  -8.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.06"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -11.49%  "
$ws.Range("E33").Value = "  -9.87%  "
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.943"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -9.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.53"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.28"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "0.0₃0648"
$ws.Range("E38").Value = "  -7.51%  "
$ws.Range("E39").Value = "  -8.37%  "
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "374.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("D43").Value = "2.629.28"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("E44").Value = "  -8.83%  "
$ws.Range("E46").Value = "  -6.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("E48").Value = "  -7.68%  "
$ws.Range("E49").Value = "  -4.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.22"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.81%  "
